{"js": "// Update the two-digit division problems in the practice table.\n// Each \"old\u00f7n=\" string is unique in the document, so a body-wide\n// search + replace for each pair is safe and unambiguous.\n\nconst pairs = [\n  [\"19\u00f76=\", \"45\u00f76=\"],\n  [\"88\u00f72=\", \"76\u00f75=\"],\n  [\"39\u00f76=\", \"90\u00f79=\"],\n  [\"61\u00f79=\", \"34\u00f77=\"],\n  [\"68\u00f75=\", \"65\u00f76=\"],\n  [\"17\u00f79=\", \"46\u00f79=\"],\n  [\"33\u00f72=\", \"90\u00f72=\"],\n  [\"18\u00f76=\", \"72\u00f75=\"],\n  [\"31\u00f72=\", \"67\u00f73=\"],\n  [\"74\u00f79=\", \"47\u00f75=\"],\n  [\"68\u00f78=\", \"28\u00f75=\"],\n  [\"33\u00f73=\", \"95\u00f77=\"],\n  [\"51\u00f74=\", \"28\u00f72=\"],\n  [\"58\u00f75=\", \"85\u00f72=\"],\n  [\"13\u00f72=\", \"40\u00f75=\"],\n  [\"87\u00f74=\", \"15\u00f73=\"],\n  [\"53\u00f79=\", \"19\u00f78=\"],\n  [\"82\u00f72=\", \"37\u00f75=\"],\n  [\"18\u00f79=\", \"89\u00f77=\"],\n  [\"45\u00f77=\", \"81\u00f78=\"],\n  [\"96\u00f78=\", \"32\u00f77=\"],\n  [\"19\u00f77=\", \"95\u00f79=\"],\n  [\"71\u00f73=\", \"85\u00f73=\"],\n  [\"88\u00f75=\", \"57\u00f78=\"],\n  [\"41\u00f78=\", \"35\u00f75=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of pairs) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the two-digit division problems in the practice table.\n# Each \"old\u00f7n=\" string is unique in the document, so a global\n# Find/Replace (wdReplaceAll) for each pair is safe and unambiguous.\n\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @{Old = \"19\u00f76=\"; New = \"45\u00f76=\"},\n    @{Old = \"88\u00f72=\"; New = \"76\u00f75=\"},\n    @{Old = \"39\u00f76=\"; New = \"90\u00f79=\"},\n    @{Old = \"61\u00f79=\"; New = \"34\u00f77=\"},\n    @{Old = \"68\u00f75=\"; New = \"65\u00f76=\"},\n    @{Old = \"17\u00f79=\"; New = \"46\u00f79=\"},\n    @{Old = \"33\u00f72=\"; New = \"90\u00f72=\"},\n    @{Old = \"18\u00f76=\"; New = \"72\u00f75=\"},\n    @{Old = \"31\u00f72=\"; New = \"67\u00f73=\"},\n    @{Old = \"74\u00f79=\"; New = \"47\u00f75=\"},\n    @{Old = \"68\u00f78=\"; New = \"28\u00f75=\"},\n    @{Old = \"33\u00f73=\"; New = \"95\u00f77=\"},\n    @{Old = \"51\u00f74=\"; New = \"28\u00f72=\"},\n    @{Old = \"58\u00f75=\"; New = \"85\u00f72=\"},\n    @{Old = \"13\u00f72=\"; New = \"40\u00f75=\"},\n    @{Old = \"87\u00f74=\"; New = \"15\u00f73=\"},\n    @{Old = \"53\u00f79=\"; New = \"19\u00f78=\"},\n    @{Old = \"82\u00f72=\"; New = \"37\u00f75=\"},\n    @{Old = \"18\u00f79=\"; New = \"89\u00f77=\"},\n    @{Old = \"45\u00f77=\"; New = \"81\u00f78=\"},\n    @{Old = \"96\u00f78=\"; New = \"32\u00f77=\"},\n    @{Old = \"19\u00f77=\"; New = \"95\u00f79=\"},\n    @{Old = \"71\u00f73=\"; New = \"85\u00f73=\"},\n    @{Old = \"88\u00f75=\"; New = \"57\u00f78=\"},\n    @{Old = \"41\u00f78=\"; New = \"35\u00f75=\"}\n)\n\nforeach ($pair in $pairs) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($pair.Old, $false, $false, $false, $false, $false, $true, 1, $false, $pair.New, 2)\n}\n"}
